$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 1044, shifting existing rows 1044-1140 down to 1046-1142.
$ws.Rows("1044:1045").Insert()

# New row 1044 (Primera quality, week of 2023-07-25)
$ws.Range("A1044").Value = 6
$ws.Range("B1044").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1044").Value = "Metropolitana"
$ws.Range("D1044").Value = 45132
$ws.Range("E1044").Value = 13
$ws.Range("F1044").Value = 100112017
$ws.Range("G1044").Value = "Apio"
$ws.Range("H1044").Value = "Americana (o)"
$ws.Range("I1044").Value = "Primera"
$ws.Range("J1044").Value = 500
$ws.Range("K1044").Value = 5000
$ws.Range("L1044").Value = 6000
$ws.Range("M1044").Value = 5540
$ws.Range("N1044").Value = "$/docena de matas"
$ws.Range("O1044").Value = "Región de Coquimbo"
$ws.Range("P1044").Value = 923
$ws.Range("Q1044").Value = 6
$ws.Range("R1044").Value = "Hortaliza"

# New row 1045 (Segunda quality, same week)
$ws.Range("A1045").Value = 6
$ws.Range("B1045").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1045").Value = "Metropolitana"
$ws.Range("D1045").Value = 45132
$ws.Range("E1045").Value = 13
$ws.Range("F1045").Value = 100112017
$ws.Range("G1045").Value = "Apio"
$ws.Range("H1045").Value = "Americana (o)"
$ws.Range("I1045").Value = "Segunda"
$ws.Range("J1045").Value = 700
$ws.Range("K1045").Value = 4000
$ws.Range("L1045").Value = 5000
$ws.Range("M1045").Value = 4500
$ws.Range("N1045").Value = "$/docena de matas"
$ws.Range("O1045").Value = "Región de Coquimbo"
$ws.Range("P1045").Value = 750
$ws.Range("Q1045").Value = 6
$ws.Range("R1045").Value = "Hortaliza"
